$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 16.06840000000001
$ws.Range("E10").Value = 16.4769
$ws.Range("E12").Value = 18.07600000000001
$ws.Range("E18").Value = 17.65350000000002
$ws.Range("E25").Value = 17.22099999999999
$ws.Range("E37").Value = 16.73240000000001
$ws.Range("E55").Value = 16.6212
$ws.Range("E68").Value = 17.47200000000002
$ws.Range("E77").Value = 18.40310000000002
$ws.Range("E78").Value = 16.70960000000003
$ws.Range("E79").Value = 18.63600000000003
$ws.Range("E80").Value = 16.69060000000002
$ws.Range("E81").Value = 16.62289999999998
$ws.Range("E82").Value = 16.90380000000001
$ws.Range("E84").Value = 16.73369999999999
$ws.Range("E101").Value = 16.89260000000002
$ws.Range("E102").Value = 16.7593
